$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$cell.NumberFormat = "@"
$cell.Value = '64.955.98'
$cell.ClearFormats()
$cell = $ws.Range('E2')
$cell.NumberFormat = "@"
$cell.Value = '  -0.48%  '
$cell.ClearFormats()
$cell = $ws.Range('D3')
$cell.NumberFormat = "@"
$cell.Value = '3.568.10'
$cell.ClearFormats()
$cell = $ws.Range('E3')
$cell.NumberFormat = "@"
$cell.Value = '  +2.40%  '
$cell.ClearFormats()
$cell = $ws.Range('E4')
$cell.NumberFormat = "@"
$cell.Value = '  +0.00%  '
$cell.ClearFormats()
$cell = $ws.Range('D5')
$cell.NumberFormat = "@"
$cell.Value = '600.18'
$cell.ClearFormats()
$cell = $ws.Range('E5')
$cell.NumberFormat = "@"
$cell.Value = '  +1.88%  '
$cell.ClearFormats()
$cell = $ws.Range('D6')
$cell.NumberFormat = "@"
$cell.Value = '135.65'
$cell.ClearFormats()
$cell = $ws.Range('E6')
$cell.NumberFormat = "@"
$cell.Value = '  -1.27%  '
$cell.ClearFormats()
$cell = $ws.Range('D7')
$cell.NumberFormat = "@"
$cell.Value = '3.566.61'
$cell.ClearFormats()
$cell = $ws.Range('E7')
$cell.NumberFormat = "@"
$cell.Value = '  +2.39%  '
$cell.ClearFormats()
$cell = $ws.Range('E8')
$cell.NumberFormat = "@"
$cell.Value = '  -0.04%  '
$cell.ClearFormats()
$cell = $ws.Range('E9')
$cell.NumberFormat = "@"
$cell.Value = '  +0.70%  '
$cell.ClearFormats()
$cell = $ws.Range('E10')
$cell.NumberFormat = "@"
$cell.Value = '  +0.63%  '
$cell.ClearFormats()
$cell = $ws.Range('E11')
$cell.NumberFormat = "@"
$cell.Value = '  -3.26%  '
$cell.ClearFormats()
$cell = $ws.Range('E12')
$cell.NumberFormat = "@"
$cell.Value = '  +0.77%  '
$cell.ClearFormats()
$cell = $ws.Range('D13')
$cell.NumberFormat = "@"
$cell.Value = '4.174.84'
$cell.ClearFormats()
$cell = $ws.Range('E13')
$cell.NumberFormat = "@"
$cell.Value = '  +2.41%  '
$cell.ClearFormats()
$cell = $ws.Range('E14')
$cell.NumberFormat = "@"
$cell.Value = '  +0.05%  '
$cell.ClearFormats()
$cell = $ws.Range('D15')
$cell.NumberFormat = "@"
$cell.Value = '3.569.67'
$cell.ClearFormats()
$cell = $ws.Range('E15')
$cell.NumberFormat = "@"
$cell.Value = '  +2.09%  '
$cell.ClearFormats()
$cell = $ws.Range('D16')
$cell.NumberFormat = "@"
$cell.Value = '27.06'
$cell.ClearFormats()
$cell = $ws.Range('E16')
$cell.NumberFormat = "@"
$cell.Value = '  +2.02%  '
$cell.ClearFormats()
$cell = $ws.Range('E17')
$cell.NumberFormat = "@"
$cell.Value = '  +0.42%  '
$cell.ClearFormats()
$cell = $ws.Range('D18')
$cell.NumberFormat = "@"
$cell.Value = '65.082.16'
$cell.ClearFormats()
$cell = $ws.Range('E18')
$cell.NumberFormat = "@"
$cell.Value = '  -0.10%  '
$cell.ClearFormats()
$cell = $ws.Range('D19')
$cell.NumberFormat = "@"
$cell.Value = '10.02'
$cell.ClearFormats()
$cell = $ws.Range('E19')
$cell.NumberFormat = "@"
$cell.Value = '  +3.10%  '
$cell.ClearFormats()
$cell = $ws.Range('D20')
$cell.NumberFormat = "@"
$cell.Value = '14.39'
$cell.ClearFormats()
$cell = $ws.Range('E20')
$cell.NumberFormat = "@"
$cell.Value = '  +3.53%  '
$cell.ClearFormats()
$cell = $ws.Range('E21')
$cell.NumberFormat = "@"
$cell.Value = '  +0.99%  '
$cell.ClearFormats()
$cell = $ws.Range('D22')
$cell.NumberFormat = "@"
$cell.Value = '389.39'
$cell.ClearFormats()
$cell = $ws.Range('E22')
$cell.NumberFormat = "@"
$cell.Value = '  +0.01%  '
$cell.ClearFormats()
$cell = $ws.Range('E23')
$cell.NumberFormat = "@"
$cell.Value = '  +4.62%  '
$cell.ClearFormats()
$cell = $ws.Range('D24')
$cell.NumberFormat = "@"
$cell.Value = '3.712.85'
$cell.ClearFormats()
$cell = $ws.Range('E24')
$cell.NumberFormat = "@"
$cell.Value = '  +2.44%  '
$cell.ClearFormats()
$cell = $ws.Range('D25')
$cell.NumberFormat = "@"
$cell.Value = '74.13'
$cell.ClearFormats()
$cell = $ws.Range('E25')
$cell.NumberFormat = "@"
$cell.Value = '  +2.11%  '
$cell.ClearFormats()
$cell = $ws.Range('E27')
$cell.NumberFormat = "@"
$cell.Value = '  +5.82%  '
$cell.ClearFormats()
$cell = $ws.Range('D28')
$cell.NumberFormat = "@"
$cell.Value = '7.74'
$cell.ClearFormats()
$cell = $ws.Range('E28')
$cell.NumberFormat = "@"
$cell.Value = '  +5.95%  '
$cell.ClearFormats()
$cell = $ws.Range('E29')
$cell.NumberFormat = "@"
$cell.Value = '  +0.12%  '
$cell.ClearFormats()
$cell = $ws.Range('E30')
$cell.NumberFormat = "@"
$cell.Value = '  +3.25%  '
$cell.ClearFormats()
$cell = $ws.Range('D31')
$cell.NumberFormat = "@"
$cell.Value = '8.45'
$cell.ClearFormats()
$cell = $ws.Range('E31')
$cell.NumberFormat = "@"
$cell.Value = '  +2.75%  '
$cell.ClearFormats()
$cell = $ws.Range('D32')
$cell.NumberFormat = "@"
$cell.Value = '1.49'
$cell.ClearFormats()
$cell = $ws.Range('E32')
$cell.NumberFormat = "@"
$cell.Value = '  +24.78%  '
$cell.ClearFormats()
$cell = $ws.Range('D33')
$cell.NumberFormat = "@"
$cell.Value = '3.569.23'
$cell.ClearFormats()
$cell = $ws.Range('E33')
$cell.NumberFormat = "@"
$cell.Value = '  +1.88%  '
$cell.ClearFormats()
$cell = $ws.Range('D34')
$cell.NumberFormat = "@"
$cell.Value = '24.05'
$cell.ClearFormats()
$cell = $ws.Range('E34')
$cell.NumberFormat = "@"
$cell.Value = '  +4.03%  '
$cell.ClearFormats()
$cell = $ws.Range('E36')
$cell.NumberFormat = "@"
$cell.Value = '  +0.49%  '
$cell.ClearFormats()
$cell = $ws.Range('E37')
$cell.NumberFormat = "@"
$cell.Value = '  +1.69%  '
$cell.ClearFormats()
$cell = $ws.Range('D38')
$cell.NumberFormat = "@"
$cell.Value = '169.29'
$cell.ClearFormats()
$cell = $ws.Range('E38')
$cell.NumberFormat = "@"
$cell.Value = '  -1.61%  '
$cell.ClearFormats()
$cell = $ws.Range('D39')
$cell.NumberFormat = "@"
$cell.Value = '1.55'
$cell.ClearFormats()
$cell = $ws.Range('E39')
$cell.NumberFormat = "@"
$cell.Value = '  +5.22%  '
$cell.ClearFormats()
$cell = $ws.Range('E40')
$cell.NumberFormat = "@"
$cell.Value = '  +5.56%  '
$cell.ClearFormats()
$cell = $ws.Range('D41')
$cell.NumberFormat = "@"
$cell.Value = '0.0808'
$cell.ClearFormats()
$cell = $ws.Range('E41')
$cell.NumberFormat = "@"
$cell.Value = '  +3.70%  '
$cell.ClearFormats()
$cell = $ws.Range('D42')
$cell.NumberFormat = "@"
$cell.Value = '27.26'
$cell.ClearFormats()
$cell = $ws.Range('E42')
$cell.NumberFormat = "@"
$cell.Value = '  +8.18%  '
$cell.ClearFormats()
$cell = $ws.Range('D43')
$cell.NumberFormat = "@"
$cell.Value = '0.826'
$cell.ClearFormats()
$cell = $ws.Range('E43')
$cell.NumberFormat = "@"
$cell.Value = '  +1.86%  '
$cell.ClearFormats()
$cell = $ws.Range('D44')
$cell.NumberFormat = "@"
$cell.Value = '42.74'
$cell.ClearFormats()
$cell = $ws.Range('E44')
$cell.NumberFormat = "@"
$cell.Value = '  +0.42%  '
$cell.ClearFormats()
$cell = $ws.Range('E45')
$cell.NumberFormat = "@"
$cell.Value = '  +0.02%  '
$cell.ClearFormats()
$cell = $ws.Range('E46')
$cell.NumberFormat = "@"
$cell.Value = '  +2.61%  '
$cell.ClearFormats()
$cell = $ws.Range('E47')
$cell.NumberFormat = "@"
$cell.Value = '  +4.21%  '
$cell.ClearFormats()
$cell = $ws.Range('D48')
$cell.NumberFormat = "@"
$cell.Value = '1.65'
$cell.ClearFormats()
$cell = $ws.Range('E48')
$cell.NumberFormat = "@"
$cell.Value = '  +1.46%  '
$cell.ClearFormats()
$cell = $ws.Range('D49')
$cell.NumberFormat = "@"
$cell.Value = '2.479.56'
$cell.ClearFormats()
$cell = $ws.Range('E49')
$cell.NumberFormat = "@"
$cell.Value = '  +11.55%  '
$cell.ClearFormats()
$cell = $ws.Range('E50')
$cell.NumberFormat = "@"
$cell.Value = '  +3.20%  '
$cell.ClearFormats()
$cell = $ws.Range('D51')
$cell.NumberFormat = "@"
$cell.Value = '2.38'
$cell.ClearFormats()
$cell = $ws.Range('E51')
$cell.NumberFormat = "@"
$cell.Value = '  +9.17%  '
$cell.ClearFormats()
